$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 (quarterly revision bumped a few totals) ---
$ws.Cells.Item(74, 2).Value = 24767   # B74
$ws.Cells.Item(74, 7).Value = 16536   # G74
$ws.Cells.Item(74, 8).Value = 10947   # H74

# --- Append new row 75 (01-04-2021 quarter) ---
# Column A holds a text label that looks like a date ("01-04-2021"). Excel's
# live-entry parser would normally convert a date-shaped string typed into a
# General-formatted cell into a date serial. Format the cell as Text first so
# the literal label is preserved (matches the rest of column A, which stores
# these labels as plain text), then drop the cell back to the default/general
# style so no extra formatting is left behind on the cell itself.
$ws.Cells.Item(75, 1).NumberFormat = "@"
$ws.Cells.Item(75, 1).Value = "01-04-2021"
$ws.Cells.Item(75, 1).ClearFormats()

$ws.Cells.Item(75, 2).Value = 24276
$ws.Cells.Item(75, 3).Value = 8446
$ws.Cells.Item(75, 4).Value = 1044
$ws.Cells.Item(75, 5).Value = 4797
$ws.Cells.Item(75, 6).Value = 2605
$ws.Cells.Item(75, 7).Value = 15830
$ws.Cells.Item(75, 8).Value = 11284
$ws.Cells.Item(75, 9).Value = 4546
